# Add a new "2022-Q3" quarter sheet to the workbook, insert its summary row
# into the "总计" sheet, and shift all other quarter sheets along (their
# underlying data is unchanged; only their sheet position moves by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计" (i.e. before
#    the sheet that is currently in position 2, "2022-Q2").
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# Header row (bold, centered, thin border - matches the other quarter sheets)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$hdr = $q3.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# Row index column (A) style used throughout these sheets
$idxStyle = $q3.Range("A2:A3")
$idxStyle.Font.Bold = $true
$idxStyle.HorizontalAlignment = -4108
$idxStyle.VerticalAlignment = -4160
$idxStyle.Borders.LineStyle = 1
$idxStyle.Borders.Weight = 2

# Columns B and D:G hold numeric-looking values stored as text in this
# workbook (fund codes keep leading zeros, percentages/sizes are text), so
# force text format before writing them.
$q3.Range("B2:B3").NumberFormat = "@"
$q3.Range("D2:G3").NumberFormat = "@"

# Row 2 - 320022 诺安研究精选股票
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "320022"
$q3.Range("C2").Value = "诺安研究精选股票"
$q3.Range("D2").Value = "6.17"
$q3.Range("E2").Value = "92.67"
$q3.Range("F2").Value = "2.82"
$q3.Range("G2").Value = "0.1740"
$q3.Range("H2").Value = 4

# Row 3 - 002291 诺安安鑫灵活配置混合
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "002291"
$q3.Range("C3").Value = "诺安安鑫灵活配置混合"
$q3.Range("D3").Value = "2.66"
$q3.Range("E3").Value = "77.38"
$q3.Range("F3").Value = "2.38"
$q3.Range("G3").Value = "0.0633"
$q3.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2. Insert the corresponding summary row into the "总计" sheet and
#    renumber the index column so it keeps counting 0,1,2,...
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
# Re-apply the original row's formatting to the freshly inserted row.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.24

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
